$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns are stored as plain text (e.g. "52.271.42",
# "340.76", "  +6.02%  "), even though many of the values look numeric.
# Assigning through .Value / .Value2 / .Formula lets Excel "smart-type"
# anything that parses as a number into a real numeric cell, which would
# flip the stored cell type away from text. Temporarily forcing the
# Price/Volume data range to a text number format keeps every write a
# literal string; ClearFormats() afterwards removes the temporary
# formatting again so no stray style survives the edit.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value2 = "52.318.55"
$ws.Range("E2").Value2 = "  +6.00%  "

# Row 3 - Ethereum
$ws.Range("D3").Value2 = "2.792.64"
$ws.Range("E3").Value2 = "  +6.37%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value2 = "  +0.06%  "

# Row 5 - Solana
$ws.Range("D5").Value2 = "116.77"
$ws.Range("E5").Value2 = "  +4.60%  "

# Row 6 - BNB
$ws.Range("D6").Value2 = "340.86"

# Row 7 - XRP
$ws.Range("D7").Value2 = "0.553"
$ws.Range("E7").Value2 = "  +5.51%  "

# Row 8 - USDC
$ws.Range("D8").Value2 = "1.00"
$ws.Range("E8").Value2 = "  +0.07%  "

# Row 9 - Cardano
$ws.Range("E9").Value2 = "  +6.04%  "

# Row 10 - Avalanche
$ws.Range("D10").Value2 = "42.06"
$ws.Range("E10").Value2 = "  +6.79%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value2 = "  +7.09%  "

# Row 12 - Chainlink
$ws.Range("D12").Value2 = "20.10"
$ws.Range("E12").Value2 = "  +0.70%  "

# Row 13 - TRON
$ws.Range("E13").Value2 = "  +2.51%  "

# Row 14 - Polkadot
$ws.Range("E14").Value2 = "  +1.26%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value2 = "3.233.99"
$ws.Range("E15").Value2 = "  +6.55%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value2 = "2.795.62"
$ws.Range("E16").Value2 = "  +6.41%  "

# Row 17 - Polygon
$ws.Range("D17").Value2 = "0.886"
$ws.Range("E17").Value2 = "  +4.06%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value2 = "52.125.09"
$ws.Range("E18").Value2 = "  +5.74%  "

# Row 19 - ImmutableX
$ws.Range("D19").Value2 = "3.21"
$ws.Range("E19").Value2 = "  +10.62%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("E20").Value2 = "  +0.94%  "

# Row 21 - Uniswap
$ws.Range("D21").Value2 = "6.97"
$ws.Range("E21").Value2 = "  +4.69%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value2 = "0.0₃0984"
$ws.Range("E22").Value2 = "  +4.06%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value2 = "278.85"
$ws.Range("E23").Value2 = "  +4.15%  "

# Row 24 - Litecoin
$ws.Range("D24").Value2 = "70.29"
$ws.Range("E24").Value2 = "  +1.94%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value2 = "2.77"
$ws.Range("E25").Value2 = "  +8.50%  "

# Row 26 - EthereumClassic
$ws.Range("E26").Value2 = "  +3.54%  "

# Row 27 - Dai
$ws.Range("E27").Value2 = "  -0.02%  "

# Row 28 - Cosmos
$ws.Range("E28").Value2 = "  +1.05%  "

# Row 29 - Toncoin
$ws.Range("D29").Value2 = "2.23"
$ws.Range("E29").Value2 = "  +1.41%  "

# Row 30 - Kaspa
$ws.Range("E30").Value2 = "  +3.56%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").Value2 = "34.86"
$ws.Range("E31").Value2 = "  +1.02%  "

# Row 32 - OKB
$ws.Range("E32").Value2 = "  +1.80%  "

# Row 33 - Filecoin
$ws.Range("D33").Value2 = "5.77"
$ws.Range("E33").Value2 = "  +5.55%  "

# Row 34 - Hedera
$ws.Range("D34").Value2 = "0.0830"
$ws.Range("E34").Value2 = "  +2.98%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value2 = "  +5.14%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value2 = "  -0.06%  "

# Row 37 - Celestia
$ws.Range("D37").Value2 = "18.99"
$ws.Range("E37").Value2 = "  +0.07%  "

# Row 38 - RenderToken
$ws.Range("D38").Value2 = "4.97"
$ws.Range("E38").Value2 = "  +0.42%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value2 = "  +5.00%  "

# Row 40 - Stacks
$ws.Range("E40").Value2 = "  +29.03%  "

# Row 41 - VeChain
$ws.Range("D41").Value2 = "0.0378"
$ws.Range("E41").Value2 = "  +13.77%  "

# Row 42 - was Stellar, now EnergySwap
$ws.Range("B42").Value2 = "EnergySwap"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value2 = "23.51"
$ws.Range("E42").Value2 = "  +3.99%  "

# Row 43 - was EnergySwap, now Stellar
$ws.Range("B43").Value2 = "Stellar"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value2 = "0.116"
$ws.Range("E43").Value2 = "  +4.45%  "

# Row 44 - WEMIXToken
$ws.Range("E44").Value2 = "  +1.35%  "

# Row 45 - Monero
$ws.Range("D45").Value2 = "124.78"
$ws.Range("E45").Value2 = "  -3.69%  "

# Row 46 - Maker
$ws.Range("D46").Value2 = "2.100.64"
$ws.Range("E46").Value2 = "  +2.08%  "

# Row 47 - NEARProtocol
$ws.Range("D47").Value2 = "3.33"
$ws.Range("E47").Value2 = "  +2.09%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").Value2 = "  +3.44%  "

# Row 49 - THORChain
$ws.Range("D49").Value2 = "5.59"
$ws.Range("E49").Value2 = "  +7.71%  "

# Row 50 - SEI
$ws.Range("D50").Value2 = "0.902"
$ws.Range("E50").Value2 = "  +21.85%  "

# Row 51 - FraxShare
$ws.Range("E51").Value2 = "  +1.74%  "

# Drop the temporary text format so the cells end up with no explicit
# style, matching the untouched rows/columns in this workbook.
$dataRange.ClearFormats()
